# Applies the "add functionality for rsm" update to Call_CBU data:
#  - CBU33 is removed from the CBU row sequence and CBU46 is introduced,
#    which shifts the row labels in column A for rows 14-21 accordingly.
#  - All call-volume figures (columns B:J, rows 2-38) are refreshed with
#    the latest reported counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the CBU row labels affected by the CBU33 -> CBU46 re-sequencing ---
$ws.Range("A14").Value = "CBU34"
$ws.Range("A15").Value = "CBU35"
$ws.Range("A16").Value = "CBU36"
$ws.Range("A17").Value = "CBU41"
$ws.Range("A18").Value = "CBU42"
$ws.Range("A19").Value = "CBU44"
$ws.Range("A20").Value = "CBU45"
$ws.Range("A21").Value = "CBU46"

# --- Refresh the call-volume figures for every data row (columns B:J) ---
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

$data = @(
    @(8, 5, 5, 4, 1, 9, 10, 1, 2),
    @(3, 3, 4, 5, 11, 12, 8, 1, 4),
    @(3, 4, 6, 2, 4, 5, 3, 3, 7),
    @(3, 2, 2, 2, 6, 8, 10, 2, 3),
    @(0, 4, 2, 1, 0, 10, 4, 9, 7),
    @(8, 6, 8, 3, 12, 6, 3, 0, 4),
    @(6, 10, 0, 5, 10, 7, 5, 0, 3),
    @(6, 9, 8, 3, 2, 5, 4, 7, 5),
    @(7, 7, 7, 0, 3, 13, 7, 5, 8),
    @(5, 4, 4, 5, 8, 6, 7, 5, 2),
    @(12, 4, 11, 0, 18, 19, 17, 9, 3),
    @(7, 2, 2, 4, 10, 14, 7, 5, 1),
    @(8, 14, 0, 2, 12, 20, 12, 0, 2),
    @(2, 9, 6, 3, 7, 9, 8, 2, 0),
    @(12, 1, 11, 3, 4, 5, 7, 5, 9),
    @(7, 8, 3, 3, 10, 9, 2, 2, 4),
    @(3, 3, 3, 3, 9, 6, 2, 5, 3),
    @(6, 3, 6, 1, 4, 4, 5, 2, 9),
    @(6, 2, 4, 14, 12, 4, 2, 0, 0),
    @(6, 9, 7, 16, 16, 11, 12, 1, 5),
    @(7, 0, 8, 0, 1, 7, 4, 0, 2),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(2, 2, 4, 0, 3, 5, 4, 2, 5),
    @(11, 4, 8, 4, 11, 3, 9, 0, 0),
    @(1, 1, 1, 1, 2, 3, 1, 0, 1),
    @(3, 2, 3, 5, 5, 9, 0, 3, 3),
    @(4, 2, 1, 0, 4, 12, 2, 3, 3),
    @(6, 4, 4, 1, 6, 7, 2, 0, 6),
    @(4, 2, 1, 4, 3, 6, 1, 3, 1),
    @(12, 15, 11, 4, 13, 21, 5, 7, 4),
    @(6, 11, 9, 8, 2, 7, 13, 10, 6),
    @(10, 7, 4, 4, 9, 11, 11, 0, 4),
    @(14, 13, 7, 3, 5, 6, 10, 1, 4),
    @(11, 6, 5, 4, 0, 4, 5, 0, 3),
    @(5, 9, 10, 8, 5, 13, 12, 0, 4),
    @(8, 7, 1, 6, 7, 5, 3, 3, 0),
    @(7, 2, 4, 4, 4, 10, 4, 1, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range("$($cols[$c])$rowNum").Value = $rowVals[$c]
    }
}
